# Rebuild the full "Pais" data block (rows 4-219) with the updated
# COVID-19 stats and the new sort order (descending by Casos totales),
# matching the refreshed export used by the site on 19-Jun-2020 15:44.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Datos actualizados a 19 de Junio de 2020 a las 15:44"

$rows = New-Object 'object[]' 216
$rows[0] = @("Estados Unidos", 2264862, 1211, 931150, 1213005, 0, 19, 120707)
$rows[1] = @("Brasil", 984315, 956, 520360, 416058, 0, 28, 47897)
$rows[2] = @("Rusia", 569063, 7972, 324406, 236816, 0, 181, 7841)
$rows[3] = @("India", 382497, 1406, 205844, 164037, 0, 12, 12616)
$rows[4] = @("Reino Unido", 300469, 0, 0, 0, 0, 0, 42288)
$rows[5] = @("España", 292348, 0, 0, 0, 0, 0, 27136)
$rows[6] = @("Peru", 244388, 0, 131190, 105737, 0, 0, 7461)
$rows[7] = @("Italia", 238159, 0, 180544, 23101, 0, 0, 34514)
$rows[8] = @("Chile", 225103, 0, 186441, 34821, 0, 0, 3841)
$rows[9] = @("Iran", 200262, 2615, 159192, 31678, 0, 120, 9392)
$rows[10] = @("Alemania", 190264, 138, 174100, 7212, 0, 6, 8952)
$rows[11] = @("Turquia", 184031, 0, 156022, 23127, 0, 0, 4882)
$rows[12] = @("Mexico", 165455, 5662, 123095, 22613, 0, 667, 19747)
$rows[13] = @("Pakistan", 165062, 4944, 61383, 100450, 0, 136, 3229)
$rows[14] = @("Francia", 158641, 0, 73887, 55151, 0, 0, 29603)
$rows[15] = @("Arabia Saudita", 150292, 4301, 95764, 53344, 0, 45, 1184)
$rows[16] = @("Banglades", 105535, 3243, 42945, 61202, 0, 45, 1388)
$rows[17] = @("Canada", 100220, 0, 62496, 29424, 0, 0, 8300)
$rows[18] = @("Catar", 85462, 1021, 65409, 19960, 0, 7, 93)
$rows[19] = @("Sudafrica", 83890, 0, 44920, 37233, 0, 0, 1737)
$rows[20] = @("China", 83325, 32, 78398, 293, 0, 0, 4634)
$rows[21] = @("Belgica", 60476, 128, 16751, 34030, 0, 12, 9695)
$rows[22] = @("Colombia", 60217, 0, 22680, 35587, 0, 0, 1950)
$rows[23] = @("Bielorrusia", 57333, 676, 35275, 21721, 0, 6, 337)
$rows[24] = @("Suecia", 56043, 0, 0, 0, 0, 0, 5053)
$rows[25] = @("Egipto", 50437, 0, 13528, 34971, 0, 0, 1938)
$rows[26] = @("Paises Bajos", 49426, 107, 0, 0, 0, 3, 6081)
$rows[27] = @("Ecuador", 49097, 0, 24123, 20887, 0, 0, 4087)
$rows[28] = @("Emiratos Arabes Unidos", 44145, 393, 30996, 12849, 0, 2, 300)
$rows[29] = @("Indonesia", 43803, 1041, 17349, 24081, 0, 34, 2373)
$rows[30] = @("Singapur", 41615, 142, 32712, 8877, 0, 0, 26)
$rows[31] = @("Kuwait", 38678, 604, 30190, 8175, 0, 5, 313)
$rows[32] = @("Portugal", 38089, 0, 24010, 12555, 0, 0, 1524)
$rows[33] = @("Argentina", 37510, 0, 10721, 25841, 0, 0, 948)
$rows[34] = @("Ucrania", 34984, 921, 16033, 17966, 0, 19, 985)
$rows[35] = @("Polonia", 31316, 301, 15698, 14284, 0, 18, 1334)
$rows[36] = @("Suiza", 31217, 17, 28900, 361, 0, 0, 1956)
$rows[37] = @("Filipinas", 28459, 660, 7378, 19951, 0, 14, 1130)
$rows[38] = @("Afganistan", 27878, 346, 7962, 19368, 0, 2, 548)
$rows[39] = @("Oman", 27670, 852, 13974, 13571, 0, 6, 125)
$rows[40] = @("Irak", 25717, 0, 11333, 13528, 0, 0, 856)
$rows[41] = @("Irlanda", 25355, 0, 22698, 943, 0, 0, 1714)
$rows[42] = @("Republica Dominicana", 24645, 0, 14293, 9717, 0, 0, 635)
$rows[43] = @("Rumania", 23400, 320, 16555, 5361, 0, 11, 1484)
$rows[44] = @("Panama", 23351, 0, 13782, 9094, 0, 0, 475)
$rows[45] = @("Bolivia", 21499, 814, 4002, 16800, 0, 18, 697)
$rows[46] = @("Barein", 20430, 0, 14696, 5677, 0, 2, 57)
$rows[47] = @("Israel", 20243, 207, 15567, 4372, 0, 1, 304)
$rows[48] = @("Armenia", 19157, 459, 8266, 10572, 0, 10, 319)
$rows[49] = @("Nigeria", 18480, 0, 6307, 11698, 0, 0, 475)
$rows[50] = @("Japon", 17668, 0, 15930, 803, 0, 0, 935)
$rows[51] = @("Austria", 17271, 48, 16141, 442, 0, 0, 688)
$rows[52] = @("Kazajistan", 16351, 474, 10195, 6051, 0, 5, 105)
$rows[53] = @("Moldavia", 13106, 0, 7525, 5132, 0, 5, 449)
$rows[54] = @("Ghana", 12929, 0, 4468, 8395, 0, 0, 66)
$rows[55] = @("Serbia", 12709, 93, 11511, 939, 0, 1, 259)
$rows[56] = @("Dinamarca", 12391, 47, 11282, 509, 0, 0, 600)
$rows[57] = @("Corea del Sur", 12306, 49, 10835, 1191, 0, 0, 280)
$rows[58] = @("Guatemala", 11868, 617, 2290, 9129, 0, 17, 449)
$rows[59] = @("Argelia", 11385, 0, 8078, 2496, 0, 0, 811)
$rows[60] = @("Azerbaiyan", 11329, 0, 6192, 4998, 0, 0, 139)
$rows[61] = @("Honduras", 10739, 440, 1179, 9217, 0, 7, 343)
$rows[62] = @("Camerun", 10638, 774, 7548, 2808, 0, 6, 282)
$rows[63] = @("Chequia", 10283, 3, 7446, 2503, 0, 0, 334)
$rows[64] = @("Marruecos", 9280, 206, 8081, 986, 0, 0, 213)
$rows[65] = @("Noruega", 8721, 13, 8138, 339, 0, 0, 244)
$rows[66] = @("Malasia", 8535, 6, 8070, 344, 0, 0, 121)
$rows[67] = @("Nepal", 8274, 426, 1402, 6850, 0, 0, 22)
$rows[68] = @("Sudan", 8020, 0, 2966, 4567, 0, 0, 487)
$rows[69] = @("Australia", 7409, 18, 6878, 429, 0, 0, 102)
$rows[70] = @("Finlandia", 7133, 14, 6200, 607, 0, 0, 326)
$rows[71] = @("Costa de Marfil", 6444, 0, 2863, 3532, 0, 0, 49)
$rows[72] = @("Uzbekistan", 5920, 153, 4273, 1628, 0, 0, 19)
$rows[73] = @("Senegal", 5639, 164, 3788, 1772, 0, 3, 79)
$rows[74] = @("Consejo Danes para los Refugiados", 5477, 194, 719, 4636, 0, 5, 122)
$rows[75] = @("Tayikistan", 5279, 0, 3762, 1466, 0, 0, 51)
$rows[76] = @("Haiti", 4916, 228, 24, 4808, 0, 2, 84)
$rows[77] = @("Guinea", 4841, 0, 3467, 1348, 0, 0, 26)
$rows[78] = @("Republica de Macedonia", 4664, 0, 1836, 2612, 0, 0, 216)
$rows[79] = @("Republica de Yibuti", 4557, 0, 3527, 987, 0, 0, 43)
$rows[80] = @("Kenia", 4374, 117, 1550, 2705, 0, 2, 119)
$rows[81] = @("Gabon", 4340, 0, 1657, 2651, 0, 0, 32)
$rows[82] = @("El Salvador", 4329, 129, 2310, 1933, 0, 4, 86)
$rows[83] = @("Luxemburgo", 4091, 0, 3940, 41, 0, 0, 110)
$rows[84] = @("Hungria", 4081, 2, 2581, 932, 0, 0, 568)
$rows[85] = @("Etiopia", 4070, 116, 1027, 2971, 0, 7, 72)
$rows[86] = @("Bulgaria", 3674, 0, 1941, 1543, 0, 0, 190)
$rows[87] = @("Venezuela", 3483, 0, 835, 2620, 0, 0, 28)
$rows[88] = @("Bosnia y Herzegovina", 3273, 99, 2241, 863, 0, 1, 169)
$rows[89] = @("Grecia", 3227, 0, 1374, 1665, 0, 0, 188)
$rows[90] = @("Tailandia", 3146, 5, 3008, 80, 0, 0, 58)
$rows[91] = @("Kirguistan", 2789, 132, 1961, 796, 0, 1, 32)
$rows[92] = @("Somalia", 2719, 0, 724, 1907, 0, 0, 88)
$rows[93] = @("Republica de Africa Central", 2605, 0, 417, 2169, 0, 0, 19)
$rows[94] = @("Mauritania", 2424, 0, 550, 1777, 0, 0, 97)
$rows[95] = @("Mayotte", 2383, 0, 2066, 288, 0, 0, 29)
$rows[96] = @("Cuba", 2295, 0, 2020, 190, 0, 0, 85)
$rows[97] = @("Croacia", 2280, 11, 2142, 31, 0, 0, 107)
$rows[98] = @("Maldivas", 2137, 0, 1759, 370, 0, 0, 8)
$rows[99] = @("Estonia", 1979, 2, 1755, 155, 0, 0, 69)
$rows[100] = @("Sri Lanka", 1948, 2, 1446, 491, 0, 0, 11)
$rows[101] = @("Costa Rica", 1939, 0, 937, 990, 0, 0, 12)
$rows[102] = @("Mali", 1923, 17, 1217, 598, 0, 1, 108)
$rows[103] = @("Albania", 1838, 50, 1114, 682, 0, 3, 42)
$rows[104] = @("Sudan del Sur", 1830, 0, 117, 1681, 0, 0, 32)
$rows[105] = @("Nicaragua", 1823, 0, 1238, 521, 0, 0, 64)
$rows[106] = @("Islandia", 1819, 3, 1801, 8, 0, 0, 10)
$rows[107] = @("Lituania", 1792, 8, 1462, 254, 0, 0, 76)
$rows[108] = @("Guayana Francesa", 1758, 0, 784, 969, 0, 0, 5)
$rows[109] = @("Guinea Ecuatorial", 1664, 0, 515, 1117, 0, 0, 32)
$rows[110] = @("Eslovaquia", 1576, 14, 1447, 101, 0, 0, 28)
$rows[111] = @("Eslovenia", 1513, 2, 1359, 45, 0, 0, 109)
$rows[112] = @("Libano", 1510, 15, 960, 518, 0, 0, 32)
$rows[113] = @("Nueva Zelanda", 1507, 0, 1482, 3, 0, 0, 22)
$rows[114] = @("Guinea-Bisau", 1492, 0, 153, 1324, 0, 0, 15)
$rows[115] = @("Madagascar", 1443, 40, 498, 932, 0, 0, 13)
$rows[116] = @("Zambia", 1416, 0, 1144, 261, 0, 0, 11)
$rows[117] = @("Paraguay", 1330, 0, 717, 600, 0, 0, 13)
$rows[118] = @("Sierra Leona", 1272, 0, 710, 511, 0, 0, 51)
$rows[119] = @("Tunez", 1146, 14, 1014, 82, 0, 0, 50)
$rows[120] = @("Hong Kong", 1128, 3, 1074, 50, 0, 0, 4)
$rows[121] = @("Letonia", 1110, 2, 903, 177, 0, 0, 30)
$rows[122] = @("Niger", 1020, 0, 901, 52, 0, 0, 67)
$rows[123] = @("Jordania", 1001, 0, 697, 295, 0, 0, 9)
$rows[124] = @("Republica de Chipre", 985, 0, 818, 148, 0, 0, 19)
$rows[125] = @("Yemen", 909, 0, 273, 388, 0, 0, 248)
$rows[126] = @("Burkina Faso", 899, 0, 810, 36, 0, 0, 53)
$rows[127] = @("Georgia", 895, 2, 741, 140, 0, 0, 14)
$rows[128] = @("Congo", 883, 0, 391, 465, 0, 0, 27)
$rows[129] = @("Principado de Andorra", 855, 0, 792, 11, 0, 0, 52)
$rows[130] = @("Republica del Chad", 854, 0, 733, 47, 0, 0, 74)
$rows[131] = @("Uruguay", 850, 0, 814, 12, 0, 0, 24)
$rows[132] = @("Cabo Verde", 823, 0, 377, 439, 0, 0, 7)
$rows[133] = @("Uganda", 755, 14, 492, 263, 0, 0, 0)
$rows[134] = @("Crucero", 712, 0, 651, 48, 0, 0, 13)
$rows[135] = @("San Marino", 696, 0, 610, 44, 0, 0, 42)
$rows[136] = @("Santo Tome y Principe", 688, 0, 191, 485, 0, 0, 12)
$rows[137] = @("Estado de Palestina", 663, 63, 437, 223, 0, 0, 3)
$rows[138] = @("Malta", 663, 0, 613, 41, 0, 0, 9)
$rows[139] = @("Mozambique", 662, 0, 175, 483, 0, 0, 4)
$rows[140] = @("Ruanda", 646, 0, 350, 294, 0, 0, 2)
$rows[141] = @("Jamaica", 638, 12, 458, 170, 0, 0, 10)
$rows[142] = @("Benin", 597, 0, 238, 348, 0, 0, 11)
$rows[143] = @("Malaui", 592, 0, 74, 510, 0, 0, 8)
$rows[144] = @("Suazilandia", 586, 0, 267, 315, 0, 0, 4)
$rows[145] = @("Togo", 547, 0, 353, 181, 0, 0, 13)
$rows[146] = @("Liberia", 542, 0, 250, 259, 0, 0, 33)
$rows[147] = @("Libia", 510, 0, 81, 419, 0, 0, 10)
$rows[148] = @("Tanzania", 509, 0, 183, 305, 0, 0, 21)
$rows[149] = @("Reunion", 502, 0, 460, 41, 0, 0, 1)
$rows[150] = @("Zimbabue", 463, 0, 63, 396, 0, 0, 4)
$rows[151] = @("Taiwan", 446, 0, 434, 5, 0, 0, 7)
$rows[152] = @("Vietnam", 349, 7, 326, 23, 0, 0, 0)
$rows[153] = @("Montenegro", 337, 0, 315, 13, 0, 0, 9)
$rows[154] = @("Mauricio", 337, 0, 325, 2, 0, 0, 10)
$rows[155] = @("Isla de Man", 336, 0, 312, 0, 0, 0, 24)
$rows[156] = @("Birmania", 286, 23, 187, 93, 0, 0, 6)
$rows[157] = @("Surinam", 277, 0, 74, 196, 0, 0, 7)
$rows[158] = @("Martinica", 221, 0, 98, 109, 0, 0, 14)
$rows[159] = @("Comoras", 210, 0, 129, 76, 0, 0, 5)
$rows[160] = @("Mongolia", 204, 3, 132, 72, 0, 0, 0)
$rows[161] = @("Islas Caimanes", 193, 0, 141, 51, 0, 0, 1)
$rows[162] = @("Siria", 187, 0, 78, 102, 0, 0, 7)
$rows[163] = @("Islas Feroe", 187, 0, 187, 0, 0, 0, 0)
$rows[164] = @("Guyana", 183, 0, 102, 69, 0, 0, 12)
$rows[165] = @("Gibraltar", 176, 0, 176, 0, 0, 0, 0)
$rows[166] = @("Guadalupe", 171, 0, 157, 0, 0, 0, 14)
$rows[167] = @("Angola", 166, 0, 64, 94, 0, 0, 8)
$rows[168] = @("Bermudas", 144, 0, 128, 7, 0, 0, 9)
$rows[169] = @("Eritrea", 142, 0, 39, 103, 0, 0, 0)
$rows[170] = @("Brunei", 141, 0, 138, 0, 0, 0, 3)
$rows[171] = @("Camboya", 129, 0, 126, 3, 0, 0, 0)
$rows[172] = @("Trinidad yTobago", 123, 0, 109, 6, 0, 0, 8)
$rows[173] = @("Burundi", 104, 0, 75, 28, 0, 0, 1)
$rows[174] = @("Bahamas", 104, 0, 74, 19, 0, 0, 11)
$rows[175] = @("Aruba", 101, 0, 98, 0, 0, 0, 3)
$rows[176] = @("Monaco", 99, 0, 94, 1, 0, 0, 4)
$rows[177] = @("Barbados", 97, 0, 85, 5, 0, 0, 7)
$rows[178] = @("Liechtenstein", 82, 0, 55, 26, 0, 0, 1)
$rows[179] = @("Botsuana", 79, 0, 25, 53, 0, 0, 1)
$rows[180] = @("San Martin (Parte Holandesa)", 77, 0, 62, 0, 0, 0, 15)
$rows[181] = @("Butan", 67, 0, 25, 42, 0, 0, 0)
$rows[182] = @("Polinesia Francesa", 60, 0, 60, 0, 0, 0, 0)
$rows[183] = @("Macao", 45, 0, 45, 0, 0, 0, 0)
$rows[184] = @("San Martin (Parte Francesa)", 42, 0, 36, 3, 0, 0, 3)
$rows[185] = @("Namibia", 40, 1, 19, 21, 0, 0, 0)
$rows[186] = @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
$rows[187] = @("Gambia", 36, 0, 24, 11, 0, 0, 1)
$rows[188] = @("Guam", 32, 0, 0, 31, 0, 0, 1)
$rows[189] = @("San Vicente y las Granadinas", 29, 0, 26, 3, 0, 0, 0)
$rows[190] = @("Antigua y Barbuda", 26, 0, 22, 1, 0, 0, 3)
$rows[191] = @("Timor Oriental", 24, 0, 24, 0, 0, 0, 0)
$rows[192] = @("Curazao", 23, 0, 19, 3, 0, 0, 1)
$rows[193] = @("Granada", 23, 0, 23, 0, 0, 0, 0)
$rows[194] = @("Belice", 22, 0, 16, 4, 0, 0, 2)
$rows[195] = @("Nueva Caledonia", 21, 0, 21, 0, 0, 0, 0)
$rows[196] = @("Santa Lucia", 19, 0, 18, 1, 0, 0, 0)
$rows[197] = @("Laos", 19, 0, 19, 0, 0, 0, 0)
$rows[198] = @("Fiyi", 18, 0, 18, 0, 0, 0, 0)
$rows[199] = @("Dominica", 18, 0, 18, 0, 0, 0, 0)
$rows[200] = @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
$rows[201] = @("San Cristobal y Nieves", 15, 0, 15, 0, 0, 0, 0)
$rows[202] = @("Groenlandia", 13, 0, 13, 0, 0, 0, 0)
$rows[203] = @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
$rows[204] = @("Islas Turcas y Caicos", 12, 0, 11, 0, 0, 0, 1)
$rows[205] = @("Santa Sede", 12, 0, 12, 0, 0, 0, 0)
$rows[206] = @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
$rows[207] = @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
$rows[208] = @("Sahara Occidental", 9, 0, 8, 0, 0, 0, 1)
$rows[209] = @("Islas Virgenes Britanicas", 8, 0, 7, 0, 0, 0, 1)
$rows[210] = @("Papua Nueva Guinea", 8, 0, 8, 0, 0, 0, 0)
$rows[211] = @("Bonaire, San Eustaquio y Saba", 7, 0, 7, 0, 0, 0, 0)
$rows[212] = @("San Bartolome", 6, 0, 6, 0, 0, 0, 0)
$rows[213] = @("Lesoto", 4, 0, 2, 2, 0, 0, 0)
$rows[214] = @("Anguila", 3, 0, 3, 0, 0, 0, 0)
$rows[215] = @("San Pedro y Miquelon", 1, 0, 1, 0, 0, 0, 0)

$data = New-Object 'object[,]' 216,8
for ($r = 0; $r -lt 216; $r++) {
    for ($c = 0; $c -lt 8; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$ws.Range("A4:H219").Value = $data

Write-Output "updated $($rows.Count) country rows"